$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (changed) date in column C moved forward one day
# (2023-12-08 -> 2023-12-09, serial 45268 -> 45269) for every data row.
$ws.Range("C2:C27").Value = 45269

# The last log entry (row 28, "A 62339-2023") was removed from the sheet.
$ws.Rows(28).Delete()

# After the deletion, the former last data row (now row 27) drops its
# explicit custom row height and reverts to the sheet's default height.
$ws.Rows(27).EntireRow.AutoFit()
